$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell as literal text (avoids Excel auto-converting
# numeric-looking strings like "1.022" into numbers) while keeping the
# cell's style/format identical to before the write.
function Set-TextCell($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextCell "D2" "28.402.46"
Set-TextCell "E2" "  +0.50%  "

Set-TextCell "D3" "1.870.94"
Set-TextCell "E3" "  -0.74%  "

Set-TextCell "D4" "1.022"
Set-TextCell "E4" "  +0.87%  "

Set-TextCell "D5" "316.75"
Set-TextCell "E5" "  +0.26%  "

Set-TextCell "E6" "  +0.84%  "

Set-TextCell "D7" "0.5104"
Set-TextCell "E7" "  -0.67%  "

Set-TextCell "D8" "0.3944"
Set-TextCell "E8" "  +0.99%  "

Set-TextCell "D9" "0.08422"
Set-TextCell "E9" "  +0.87%  "

Set-TextCell "D10" "1.108"
Set-TextCell "E10" "  -1.27%  "

Set-TextCell "D11" "41.99"
Set-TextCell "E11" "  +0.51%  "

Set-TextCell "D12" "6.231"
Set-TextCell "E12" "  +0.04%  "

Set-TextCell "D13" "1.874.93"
Set-TextCell "E13" "  -0.82%  "

Set-TextCell "D14" "20.42"
Set-TextCell "E14" "  -0.73%  "

Set-TextCell "B15" "BinanceUSD"
Set-TextCell "C15" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextCell "D15" "1.022"
Set-TextCell "E15" "  +0.99%  "

Set-TextCell "B16" "Chainlink"
Set-TextCell "C16" "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextCell "D16" "7.217"
Set-TextCell "E16" "  -0.59%  "

Set-TextCell "D17" "0.00001109"
Set-TextCell "E17" "  +0.41%  "

Set-TextCell "D18" "90.87"
Set-TextCell "E18" "  -0.29%  "

Set-TextCell "D19" "0.06765"
Set-TextCell "E19" "  +0.96%  "

Set-TextCell "D20" "17.68"
Set-TextCell "E20" "  -0.80%  "

Set-TextCell "E21" "  +0.74%  "

Set-TextCell "D22" "5.925"
Set-TextCell "E22" "  -1.75%  "

Set-TextCell "D23" "28.471.05"
Set-TextCell "E23" "  +0.62%  "

Set-TextCell "D24" "11.14"
Set-TextCell "E24" "  +0.06%  "

Set-TextCell "D25" "2.289"
Set-TextCell "E25" "  +0.56%  "

Set-TextCell "D26" "2.092.87"
Set-TextCell "E26" "  -0.42%  "

Set-TextCell "D27" "161.42"
Set-TextCell "E27" "  +0.56%  "

Set-TextCell "D28" "20.61"
Set-TextCell "E28" "  -1.05%  "

Set-TextCell "D29" "2.338"
Set-TextCell "E29" "  -4.32%  "

Set-TextCell "D30" "126.89"
Set-TextCell "E30" "  +0.12%  "

Set-TextCell "E31" "  -0.67%  "

Set-TextCell "D32" "1.034"
Set-TextCell "E32" "  -0.33%  "

Set-TextCell "D33" "5.743"
Set-TextCell "E33" "  -2.08%  "

Set-TextCell "D34" "3.646"

Set-TextCell "D35" "0.02428"
Set-TextCell "E35" "  -0.55%  "

Set-TextCell "D36" "0.06449"
Set-TextCell "E36" "  -1.65%  "

Set-TextCell "D37" "0.2171"
Set-TextCell "E37" "  -1.76%  "

Set-TextCell "D38" "8.777"
Set-TextCell "E38" "  -7.17%  "

Set-TextCell "D39" "1.262"
Set-TextCell "E39" "  +1.39%  "

Set-TextCell "D40" "1.179"
Set-TextCell "E40" "  -1.67%  "

Set-TextCell "D41" "0.6359"
Set-TextCell "E41" "  -2.02%  "

Set-TextCell "D42" "4.974"
Set-TextCell "E42" "  -0.65%  "

Set-TextCell "D43" "11.18"
Set-TextCell "E43" "  -0.40%  "

Set-TextCell "D44" "0.6026"
Set-TextCell "E44" "  -1.12%  "

Set-TextCell "D45" "13.01"
Set-TextCell "E45" "  -0.98%  "

Set-TextCell "D46" "3.708"
Set-TextCell "E46" "  +0.26%  "

Set-TextCell "D47" "1.984"
Set-TextCell "E47" "  -1.54%  "

Set-TextCell "D48" "1.204"
Set-TextCell "E48" "  -6.09%  "

Set-TextCell "D49" "121.94"
Set-TextCell "E49" "  +0.57%  "

Set-TextCell "D50" "1.202"
Set-TextCell "E50" "  -2.90%  "

Set-TextCell "D51" "0.06842"
Set-TextCell "E51" "  -1.03%  "
